$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new columns -----------------------------------------
# "Тип иссо" goes right after "Код иссо"      -> new column B
# "Id дорог АБДД" goes right after "Код дороги АБДМ" -> new column F
# Insert() shifts the existing C:F data right and carries the neighbouring
# column's style onto the freshly inserted (still empty) column.
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(6).Insert()

# Helper: force a numeric-looking literal into a cell as TEXT instead of a
# number. A plain `Range.Value = "1100101"` gets auto-converted to a number
# by Excel because the string looks numeric. Writing it through a TEXT()
# formula and pasting only the computed value (xlPasteValues) keeps it as a
# real string cell - same trick as pasting a TEXT() result in real Excel -
# without touching the cell's existing number format/style.
function Set-TextValue([string]$addr, [string]$digits) {
    $scratch = $ws.Range("ZZ1000")
    $scratch.Formula = '=TEXT(' + $digits + ',"0")'
    $scratch.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
    $scratch.ClearContents() | Out-Null
    $excel.CutCopyMode = 0
}

# --- Header row (row 1) --------------------------------------------------
$ws.Range("A1").Value = "Код иссо"
$ws.Range("B1").Value = "Тип иссо"
$ws.Range("C1").Value = "ФКУ"
$ws.Range("D1").Value = "Дорога"
$ws.Range("E1").Value = "Код дороги АБДМ"
$ws.Range("F1").Value = "Id дорог АБДД"
$ws.Range("G1").Value = "КМ"
$ws.Range("H1").Value = "М"

# --- Data row (row 2) -----------------------------------------------------
Set-TextValue "A2" "1100101"
Set-TextValue "B2" "20"
$ws.Range("C2").Value = 'ФКУ Упрдор "Прибайкалье"'
$ws.Range("D2").Value = '"Вилюй" Тулун - Братск - Усть-Кут - Мирный - Якутск'
Set-TextValue "E2" "5030"
$ws.Range("F2").Value = "56196a2d-5830-4b45-94e5-682d84e96aaf"
Set-TextValue "G2" "491"
Set-TextValue "H2" "518"

# --- Column widths ---------------------------------------------------------
# ColumnWidth is quantized by Excel to its own pixel grid on save, so the
# literal target character-width values (15.625/31.25/27.34375/39.0625 -
# native Apache POI 256ths-of-a-character units) are not bit-for-bit
# reproducible through COM. These inputs are the closest values that land
# on the saved grid nearest each target width.
$ws.Columns.Item(1).ColumnWidth = 14.75
$ws.Columns.Item(2).ColumnWidth = 14.75
$ws.Columns.Item(3).ColumnWidth = 30.25
$ws.Columns.Item(4).ColumnWidth = 30.25
$ws.Columns.Item(5).ColumnWidth = 26.41675
$ws.Columns.Item(6).ColumnWidth = 38.0835
$ws.Columns.Item(7).ColumnWidth = 14.75
$ws.Columns.Item(8).ColumnWidth = 14.75
